# Update link to privacy policy document
# (bvq_03_consent.xlsx - "survey" sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Activate()

# --- Add real hyperlinks to the catalan & spanish cells (this sets the cell text too, ------
# so we overwrite Value again afterwards with the desired markdown text) -------------------
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.upf.edu/web/cbclab/politica-privacitat", "", "", "https://www.upf.edu/web/cbclab/politica-privacitat")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.upf.edu/web/cbclab/politica-privacitat", "", "", "https://www.upf.edu/web/cbclab/politica-privacitat")

# --- Update the three privacy-policy link texts -------------------------------------------
# Row 5 = link_catalan, Row 6 = link_spanish, Row 7 = link_english
$ws.Range("C5").Value = "**[Política de privacitat](https://www.upf.edu/web/cbclab/politica-privacitat)**"
$ws.Range("C6").Value = "**[Política de privacidad](https://www.upf.edu/web/cbclab/politica-privacitat)**"
$ws.Range("C7").Value = "**[Privacy policy](https://www.upf.edu/web/cbclab/politica-privacitat)**"

# --- Row heights are now much shorter since the links are one-liners ----------------------
$ws.Rows.Item(5).RowHeight = 19.7
$ws.Rows.Item(6).RowHeight = 19.7

# --- Update the visible selection / scroll position ----------------------------------------
$ws.Range("C8").Select()
